# Update the acquisition-timestamp column (A) for the data rows on the
# "ランサーズ" sheet from 2025-10-30 06:26:41 to 2025-10-30 06:36:34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-30 06:26:41"
$newTimestamp = "2025-10-30 06:36:34"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
